$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 275.92307
$ws.Range("I9").Value = 217.09091
$ws.Range("K9").Value = 217.09091
$ws.Range("M9").Value = -48.09091000000001
# Row 18
$ws.Range("H18").Value = 3090
$ws.Range("J18").Value = 3990
$ws.Range("L18").Value = 3990
$ws.Range("N18").Value = -4558
# Row 53
$ws.Range("H53").Value = 1036
$ws.Range("I53").Value = 979.6
$ws.Range("J53").Value = 1130
$ws.Range("K53").Value = 979.6
$ws.Range("L53").Value = 1130
$ws.Range("M53").Value = -342.6
$ws.Range("N53").Value = -2404
# Row 96
$ws.Range("H96").Value = 968736.4399999999
$ws.Range("I96").Value = 1432.8334
$ws.Range("K96").Value = 4298.5002
$ws.Range("M96").Value = -2925.5002
# Row 101
$ws.Range("H101").Value = 1085.1578
$ws.Range("I101").Value = 617.7143
$ws.Range("J101").Value = 1357.8334
$ws.Range("K101").Value = 1853.1429
$ws.Range("L101").Value = 4073.5002
$ws.Range("M101").Value = -231.1428999999998
$ws.Range("N101").Value = -7317.5002
# Row 103
$ws.Range("I103").Value = 488.5
$ws.Range("J103").Value = 55557360
$ws.Range("K103").Value = 1465.5
$ws.Range("L103").Value = 166672080
$ws.Range("M103").Value = -879.5
$ws.Range("N103").Value = -166673252
# Row 132
$ws.Range("H132").Value = 2118.513
$ws.Range("I132").Value = 2152.3713
$ws.Range("K132").Value = 6457.113899999999
$ws.Range("M132").Value = -3927.113899999999
# Row 137
$ws.Range("H137").Value = 2675.7097
$ws.Range("I137").Value = 1869
$ws.Range("K137").Value = 5607
$ws.Range("M137").Value = -3057
# Row 138
$ws.Range("H138").Value = 5205.2974
$ws.Range("I138").Value = 1769.6666
$ws.Range("J138").Value = 11548
$ws.Range("K138").Value = 5308.9998
$ws.Range("L138").Value = 34644
$ws.Range("M138").Value = -168.9997999999996
$ws.Range("N138").Value = -44924

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 7323.2896
$ws.Range("I32").Value = 7087.6
$ws.Range("K32").Value = 7087.6
$ws.Range("M32").Value = -6800.6
# Row 45
$ws.Range("H45").Value = 2978728.5
$ws.Range("I45").Value = 4763257.5
$ws.Range("J45").Value = 4513.3335
$ws.Range("K45").Value = 4763257.5
$ws.Range("L45").Value = 4513.3335
$ws.Range("M45").Value = -4762880.5
$ws.Range("N45").Value = -5267.3335
# Row 61
$ws.Range("H61").Value = 8088129.5
$ws.Range("I61").Value = 10012358
$ws.Range("K61").Value = 10012358
$ws.Range("M61").Value = -10012146
# Row 74
$ws.Range("H74").Value = 1731.4073
$ws.Range("I74").Value = 1593.1818
$ws.Range("J74").Value = 2339.6
$ws.Range("K74").Value = 1593.1818
$ws.Range("L74").Value = 2339.6
$ws.Range("M74").Value = -719.1818000000001
$ws.Range("N74").Value = -4087.6
# Row 77
$ws.Range("H77").Value = 1731.4073
$ws.Range("I77").Value = 1593.1818
$ws.Range("J77").Value = 2339.6
$ws.Range("K77").Value = 7965.909000000001
$ws.Range("L77").Value = 11698
$ws.Range("M77").Value = -3597.909000000001
$ws.Range("N77").Value = -20434
# Row 97
$ws.Range("H97").Value = 1436.7097
$ws.Range("I97").Value = 817.5789
$ws.Range("K97").Value = 817.5789
$ws.Range("M97").Value = -321.5789
# Row 110
$ws.Range("H110").Value = 5254.9614
$ws.Range("I110").Value = 5015.091
$ws.Range("K110").Value = 5015.091
$ws.Range("M110").Value = -2970.091
# Row 132
$ws.Range("H132").Value = 1697245.9
$ws.Range("I132").Value = 2252.6326
$ws.Range("K132").Value = 6757.8978
$ws.Range("M132").Value = -4227.8978
# Row 136
$ws.Range("H136").Value = 8088129.5
$ws.Range("I136").Value = 10012358
$ws.Range("K136").Value = 30037074
$ws.Range("M136").Value = -30034524

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 5061.7837
$ws.Range("I20").Value = 5130.276
$ws.Range("K20").Value = 5130.276
$ws.Range("M20").Value = -4883.276
# Row 86
$ws.Range("H86").Value = 735384.7
$ws.Range("I86").Value = 918097.75
$ws.Range("K86").Value = 918097.75
$ws.Range("M86").Value = -916974.75
# Row 89
$ws.Range("H89").Value = 735384.7
$ws.Range("I89").Value = 918097.75
$ws.Range("K89").Value = 4590488.75
$ws.Range("M89").Value = -4584872.75
# Row 99
$ws.Range("H99").Value = 1740.8667
$ws.Range("I99").Value = 1023.7778
$ws.Range("J99").Value = 2816.5
$ws.Range("K99").Value = 1023.7778
$ws.Range("L99").Value = 2816.5
$ws.Range("M99").Value = 474.2222
$ws.Range("N99").Value = -5812.5
# Row 105
$ws.Range("H105").Value = 419831.47
$ws.Range("I105").Value = 539453.3
$ws.Range("K105").Value = 539453.3
$ws.Range("M105").Value = -537706.3
# Row 134
$ws.Range("H134").Value = 9092496
$ws.Range("I134").Value = 1182.625
$ws.Range("K134").Value = 3547.875
$ws.Range("M134").Value = -1012.875

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 45458384
$ws.Range("I31").Value = 62502440
$ws.Range("K31").Value = 62502440
$ws.Range("M31").Value = -62502145
# Row 34
$ws.Range("H34").Value = 45458384
$ws.Range("I34").Value = 62502440
$ws.Range("K34").Value = 62502440
$ws.Range("M34").Value = -62502238
# Row 50
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
# Row 59
$ws.Range("H59").Value = 107499.5
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()
# Row 141
$ws.Range("H141").Value = 421817.53
$ws.Range("J141").Value = 482221.44
$ws.Range("L141").Value = 482221.44
$ws.Range("N141").Value = -492581.44

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 114.76471
$ws.Range("I2").Value = 97.90909000000001
$ws.Range("J2").Value = 145.66667
$ws.Range("K2").Value = 587.4545400000001
$ws.Range("L2").Value = 874.0000200000001
$ws.Range("M2").Value = -474.4545400000001
$ws.Range("N2").Value = -1100.00002
# Row 132
$ws.Range("H132").Value = 1524.2916
$ws.Range("I132").Value = 1042
$ws.Range("J132").Value = 3357
$ws.Range("K132").Value = 9378
$ws.Range("L132").Value = 30213
$ws.Range("M132").Value = -6848
$ws.Range("N132").Value = -35273

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 5496.048
$ws.Range("I122").Value = 4162.6113
$ws.Range("J122").Value = 13496.667
$ws.Range("K122").Value = 12487.8339
$ws.Range("L122").Value = 40490.001
$ws.Range("M122").Value = -10037.8339
$ws.Range("N122").Value = -45390.001
# Row 132
$ws.Range("H132").Value = 2175589
$ws.Range("I132").Value = 1727.7906
$ws.Range("K132").Value = 5183.3718
$ws.Range("M132").Value = -2653.3718
# Row 141
$ws.Range("H141").Value = 47898.2
$ws.Range("J141").Value = 47898.2
$ws.Range("L141").Value = 47898.2
$ws.Range("N141").Value = -58258.2

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 18873480
$ws.Range("I22").Value = 26422392
$ws.Range("K22").Value = 26422392
$ws.Range("M22").Value = -26422097
# Row 27
$ws.Range("H27").Value = 18873480
$ws.Range("I27").Value = 26422392
$ws.Range("K27").Value = 26422392
$ws.Range("M27").Value = -26422285
# Row 55
$ws.Range("H55").Value = 1472.2354
$ws.Range("I55").Value = 1291
$ws.Range("J55").Value = 1599.1
$ws.Range("K55").Value = 1291
$ws.Range("L55").Value = 1599.1
$ws.Range("M55").Value = -1118
$ws.Range("N55").Value = -1945.1
# Row 57
$ws.Range("H57").Value = 34080.75
$ws.Range("I57").Value = 34080.75
$ws.Range("K57").Value = 34080.75
$ws.Range("M57").Value = -33514.75
# Row 122
$ws.Range("H122").Value = 3716.68
$ws.Range("I122").Value = 3464.5957
$ws.Range("J122").Value = 7666
$ws.Range("K122").Value = 10393.7871
$ws.Range("L122").Value = 22998
$ws.Range("M122").Value = -7943.7871
$ws.Range("N122").Value = -27898
# Row 136
$ws.Range("H136").Value = 3389.2
$ws.Range("I136").Value = 2355.75
$ws.Range("J136").Value = 5456.1
$ws.Range("K136").Value = 7067.25
$ws.Range("L136").Value = 16368.3
$ws.Range("M136").Value = -4517.25
$ws.Range("N136").Value = -21468.3

$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 52166.5
$ws.Range("J2").Value = 19000
$ws.Range("L2").Value = 19000
$ws.Range("N2").Value = -19224
# Row 113
$ws.Range("H113").Value = 629.59375
$ws.Range("I113").Value = 525.125
$ws.Range("K113").Value = 1575.375
$ws.Range("M113").Value = 594.625
# Row 122
$ws.Range("H122").Value = 2522
$ws.Range("I122").Value = 2112.8235
$ws.Range("J122").Value = 6000
$ws.Range("K122").Value = 6338.470499999999
$ws.Range("L122").Value = 18000
$ws.Range("M122").Value = -3888.470499999999
$ws.Range("N122").Value = -22900
# Row 132
$ws.Range("H132").Value = 234915.33
$ws.Range("I132").Value = 2181.4412
$ws.Range("K132").Value = 6544.323600000001
$ws.Range("M132").Value = -4014.323600000001
